$wb = $excel.ActiveWorkbook

# --- Productdata sheet: StartingInventories (C) and SetupCosts (E) columns ---
$wsProductdata = $wb.Worksheets.Item("Productdata")
$wsProductdata.Range("C2").Value = 0
$wsProductdata.Range("E2").Value = 0.3563999999999999
$wsProductdata.Range("C3").Value = 5
$wsProductdata.Range("E3").Value = 0.11655
$wsProductdata.Range("C4").Value = 5
$wsProductdata.Range("E4").Value = 0.1166888888888889
$wsProductdata.Range("C5").Value = 0
$wsProductdata.Range("E5").Value = 0.1160722222222222
$wsProductdata.Range("C6").Value = 5
$wsProductdata.Range("E6").Value = 0.1138777777777778
$wsProductdata.Range("C7").Value = 5
$wsProductdata.Range("E7").Value = 0.1141333333333333
$wsProductdata.Range("C8").Value = 5
$wsProductdata.Range("E8").Value = 0.1136444444444445
$wsProductdata.Range("C9").Value = 5
$wsProductdata.Range("E9").Value = 0.1202277777777778
$wsProductdata.Range("C10").Value = 5
$wsProductdata.Range("E10").Value = 0.12
$wsProductdata.Range("C11").Value = 5
$wsProductdata.Range("E11").Value = 0.1193833333333333
$wsProductdata.Range("C12").Value = 5
$wsProductdata.Range("E12").Value = 0.1201666666666667
$wsProductdata.Range("C13").Value = 0
$wsProductdata.Range("E13").Value = 0.3413666666666666
$wsProductdata.Range("C14").Value = 0
$wsProductdata.Range("E14").Value = 0.1140277777777778
$wsProductdata.Range("C15").Value = 0
$wsProductdata.Range("E15").Value = 0.1124944444444445
$wsProductdata.Range("C16").Value = 0
$wsProductdata.Range("E16").Value = 0.1113222222222222
$wsProductdata.Range("C17").Value = 0
$wsProductdata.Range("E17").Value = 0.11165
$wsProductdata.Range("C18").Value = 0
$wsProductdata.Range("E18").Value = 0.1119
$wsProductdata.Range("C19").Value = 0
$wsProductdata.Range("E19").Value = 0.1114444444444444
$wsProductdata.Range("C20").Value = 0
$wsProductdata.Range("E20").Value = 0.1315
$wsProductdata.Range("C21").Value = 0
$wsProductdata.Range("E21").Value = 0.1396
$wsProductdata.Range("C22").Value = 0
$wsProductdata.Range("E22").Value = 0.1728
$wsProductdata.Range("C23").Value = 0
$wsProductdata.Range("E23").Value = 0.5317

# --- Capacity sheet: column B ---
$wsCapacity = $wb.Worksheets.Item("Capacity")
$wsCapacity.Range("B2").Value = 150
$wsCapacity.Range("B3").Value = 40
$wsCapacity.Range("B4").Value = 30
$wsCapacity.Range("B5").Value = 30
$wsCapacity.Range("B6").Value = 10
$wsCapacity.Range("B7").Value = 40
$wsCapacity.Range("B8").Value = 20
$wsCapacity.Range("B9").Value = 30
$wsCapacity.Range("B10").Value = 20
$wsCapacity.Range("B11").Value = 20
$wsCapacity.Range("B12").Value = 30
$wsCapacity.Range("B13").Value = 90
$wsCapacity.Range("B14").Value = 30
$wsCapacity.Range("B15").Value = 40
$wsCapacity.Range("B16").Value = 10
$wsCapacity.Range("B17").Value = 50
$wsCapacity.Range("B18").Value = 50
$wsCapacity.Range("B19").Value = 30
$wsCapacity.Range("B20").Value = 90
$wsCapacity.Range("B21").Value = 450
$wsCapacity.Range("B22").Value = 450
$wsCapacity.Range("B23").Value = 180

# --- ProcessingTime sheet: diagonal cells ---
$wsProcessingTime = $wb.Worksheets.Item("ProcessingTime")
$wsProcessingTime.Range("B2").Value = 5
$wsProcessingTime.Range("C3").Value = 4
$wsProcessingTime.Range("E5").Value = 3
$wsProcessingTime.Range("F6").Value = 1
$wsProcessingTime.Range("I9").Value = 3
$wsProcessingTime.Range("J10").Value = 2
$wsProcessingTime.Range("K11").Value = 2
$wsProcessingTime.Range("L12").Value = 3
$wsProcessingTime.Range("M13").Value = 3
$wsProcessingTime.Range("O15").Value = 4
$wsProcessingTime.Range("P16").Value = 1
$wsProcessingTime.Range("R18").Value = 5
$wsProcessingTime.Range("S19").Value = 3
$wsProcessingTime.Range("T20").Value = 1
$wsProcessingTime.Range("V22").Value = 5
$wsProcessingTime.Range("W23").Value = 2
